# FMBL_QTR_FIN.xlsx - "Doing Updates for Financials"
# The quarterly financial tables on sheet "FMBL" gain a new quarter of data.
# A new quarter's column is inserted immediately to the left of the existing
# "Period Ending" data (old column D), pushing the prior 8 quarters (D:K)
# right by two columns (to F:M) and adding two new quarters of figures
# in the freed-up columns D and E for every line item, across all three
# statements (Income Statement, Balance Sheet, Cash Flow Statement).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FMBL")

# 1) Insert two blank columns at D:E; everything from D:K moves to F:M.
$ws.Columns("D:E").Insert()

# 2) The newly inserted columns come in blank/unformatted. Give every row in
#    D:E the same number format / font / alignment as its data (now in column
#    F, the former column D) so the new quarters look like the rest of the
#    table (date header rows get the date format, all other rows get the
#    right-aligned #,##0 numeric format).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Income Statement: newest quarter -> column D, next quarter -> column E
$ws.Range("D7").Value = 43373
$ws.Range("E7").Value = 43281
$ws.Range("D8").Value = 66600
$ws.Range("E8").Value = 63700
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 7800
$ws.Range("E17").Value = 5300
$ws.Range("D18").Value = 58800
$ws.Range("E18").Value = 58400
$ws.Range("D20").Value = -37200
$ws.Range("E20").Value = -28500
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 21600
$ws.Range("E23").Value = 29900
$ws.Range("D24").Value = 900
$ws.Range("E24").Value = 8100
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 20700
$ws.Range("E26").Value = 21800
$ws.Range("D27").Value = 20700
$ws.Range("E27").Value = 21800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 37200
$ws.Range("E32").Value = 28500
$ws.Range("D33").Value = 20700
$ws.Range("E33").Value = 21800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 20700
$ws.Range("E35").Value = 21800

# Balance Sheet: newest quarter -> column D, next quarter -> column E
$ws.Range("D38").Value = 43373
$ws.Range("E38").Value = 43281
$ws.Range("D41").Value = 59500
$ws.Range("E41").Value = 61200
$ws.Range("D42").Value = 42000
$ws.Range("E42").Value = 258700
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 99100
$ws.Range("E48").Value = 97500
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 21700
$ws.Range("E52").Value = 21300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 7309200
$ws.Range("E54").Value = 7274400
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6297600
$ws.Range("E66").Value = 6279900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 898800
$ws.Range("E72").Value = 881200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1011600
$ws.Range("E76").Value = 994500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

# Cash Flow Statement: newest quarter -> column D, next quarter -> column E
$ws.Range("D80").Value = 43373
$ws.Range("E80").Value = 43281
$ws.Range("D81").Value = 20700
$ws.Range("E81").Value = 21800
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("E89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 0

